$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column D (D:F), shifting old D:K data to G:N
$ws.Range("D:F").EntireColumn.Insert()

# Carry number/date formatting from the (now-shifted) old columns into the new D:F columns
# for every contiguous block of data rows (rows 7-35, 38-77, 80-102); rows 5,6,37,79 never
# held D:K data and must stay untouched.
$ws.Range("G7:N35").Copy()
$ws.Range("D7:F35").PasteSpecial(-4122)
$ws.Range("G38:N77").Copy()
$ws.Range("D38:F77").PasteSpecial(-4122)
$ws.Range("G80:N102").Copy()
$ws.Range("D80:F102").PasteSpecial(-4122)

# Populate the new quarters data (cols D:F) plus a handful of restated historical
# figures in the shifted columns (G:N), cell by cell, matching the refreshed financials.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("D8").Value = 3100
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = 100
$ws.Range("G8").Value = 1900
$ws.Range("H8").Value = 1300
$ws.Range("D9").Value = 3100
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("D14").Value = 400
$ws.Range("E14").Value = 500
$ws.Range("F14").Value = -3000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("D17").Value = 3700
$ws.Range("E17").Value = 900
$ws.Range("F17").Value = -2400
$ws.Range("G17").Value = 2100
$ws.Range("I17").Value = -800
$ws.Range("J17").Value = -700
$ws.Range("D18").Value = -600
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = 2400
$ws.Range("H18").Value = 1900
$ws.Range("I18").Value = 800
$ws.Range("D20").Value = 100
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = 100
$ws.Range("D21").Value = -500
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = 2500
$ws.Range("H21").Value = 1800
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 100
$ws.Range("I22").Value = 0
$ws.Range("D23").Value = -600
$ws.Range("E23").Value = -800
$ws.Range("F23").Value = 2500
$ws.Range("H23").Value = 1700
$ws.Range("I23").Value = 800
$ws.Range("D24").Value = -400
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = 900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("D26").Value = -200
$ws.Range("E26").Value = -800
$ws.Range("F26").Value = 1600
$ws.Range("I26").Value = 800
$ws.Range("D27").Value = -200
$ws.Range("E27").Value = -800
$ws.Range("F27").Value = 1600
$ws.Range("I27").Value = 800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("D32").Value = -100
$ws.Range("E32").Value = "NA"
$ws.Range("F32").Value = -100
$ws.Range("D33").Value = -200
$ws.Range("E33").Value = -800
$ws.Range("F33").Value = 1600
$ws.Range("I33").Value = 800
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("D35").Value = -200
$ws.Range("E35").Value = -800
$ws.Range("F35").Value = 1600
$ws.Range("I35").Value = 800
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("D41").Value = 3500
$ws.Range("E41").Value = 2500
$ws.Range("F41").Value = 4500
$ws.Range("G41").Value = 5500
$ws.Range("H41").Value = 4100
$ws.Range("I41").Value = 5300
$ws.Range("J41").Value = 7500
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("D43").Value = 200
$ws.Range("E43").Value = 700
$ws.Range("F43").Value = 700
$ws.Range("H43").Value = 1300
$ws.Range("D44").Value = 2900
$ws.Range("E44").Value = 400
$ws.Range("F44").Value = 400
$ws.Range("G44").Value = 300
$ws.Range("H44").Value = 9300
$ws.Range("I44").Value = 1800
$ws.Range("J44").Value = 3200
$ws.Range("D45").Value = 3800
$ws.Range("E45").Value = 8500
$ws.Range("F45").Value = 9600
$ws.Range("G45").Value = 7800
$ws.Range("D46").Value = 10400
$ws.Range("E46").Value = 12100
$ws.Range("F46").Value = 15100
$ws.Range("G46").Value = 14200
$ws.Range("H46").Value = 14800
$ws.Range("I46").Value = 7500
$ws.Range("J46").Value = 11500
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 100
$ws.Range("F47").Value = 100
$ws.Range("D48").Value = 16400
$ws.Range("E48").Value = 16400
$ws.Range("F48").Value = 13700
$ws.Range("G48").Value = 10900
$ws.Range("H48").Value = 10800
$ws.Range("I48").Value = 16500
$ws.Range("J48").Value = 11900
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("D52").Value = 200
$ws.Range("E52").Value = 200
$ws.Range("F52").Value = 200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("D54").Value = 27000
$ws.Range("E54").Value = 28800
$ws.Range("F54").Value = 29100
$ws.Range("G54").Value = 25500
$ws.Range("H54").Value = 25800
$ws.Range("I54").Value = 24300
$ws.Range("J54").Value = 23700
$ws.Range("D57").Value = 700
$ws.Range("E57").Value = 1700
$ws.Range("F57").Value = 1500
$ws.Range("H57").Value = 800
$ws.Range("I57").Value = 1500
$ws.Range("D58").Value = 1000
$ws.Range("E58").Value = 1000
$ws.Range("F58").Value = 600
$ws.Range("D59").Value = 200
$ws.Range("E59").Value = 300
$ws.Range("F59").Value = 400
$ws.Range("J59").Value = 900
$ws.Range("D60").Value = 1900
$ws.Range("E60").Value = 3000
$ws.Range("F60").Value = 2500
$ws.Range("H60").Value = 1100
$ws.Range("I60").Value = 2000
$ws.Range("J60").Value = 2200
$ws.Range("D61").Value = 1300
$ws.Range("E61").Value = 1400
$ws.Range("F61").Value = 1300
$ws.Range("H61").Value = 1800
$ws.Range("I61").Value = 1700
$ws.Range("J61").Value = 1700
$ws.Range("D62").Value = 3400
$ws.Range("E62").Value = 3800
$ws.Range("F62").Value = 3800
$ws.Range("G62").Value = 3000
$ws.Range("H62").Value = 2900
$ws.Range("I62").Value = 1900
$ws.Range("J62").Value = 1900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("D66").Value = 6700
$ws.Range("E66").Value = 8200
$ws.Range("F66").Value = 7700
$ws.Range("G66").Value = 5500
$ws.Range("H66").Value = 5800
$ws.Range("I66").Value = 5700
$ws.Range("J66").Value = 5800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("D72").Value = -6600
$ws.Range("E72").Value = -6400
$ws.Range("F72").Value = -5600
$ws.Range("G72").Value = -7000
$ws.Range("H72").Value = -7000
$ws.Range("I72").Value = -8000
$ws.Range("J72").Value = -8900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("D76").Value = 20300
$ws.Range("E76").Value = 20600
$ws.Range("F76").Value = 21400
$ws.Range("G76").Value = 20000
$ws.Range("H76").Value = 20000
$ws.Range("I76").Value = 18700
$ws.Range("J76").Value = 17800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("D81").Value = -200
$ws.Range("E81").Value = -800
$ws.Range("F81").Value = 1600
$ws.Range("I81").Value = 800
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("F83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("D89").Value = 1100
$ws.Range("E89").Value = -1600
$ws.Range("F89").Value = -800
$ws.Range("I89").Value = -1700
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = -800
$ws.Range("F91").Value = -200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = -800
$ws.Range("F94").Value = -300
$ws.Range("I94").Value = -400
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("F96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("D100").Value = -100
$ws.Range("E100").Value = 400
$ws.Range("F100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("D102").Value = 1000
$ws.Range("E102").Value = -2000
$ws.Range("F102").Value = -1100
$ws.Range("G102").Value = 1400
$ws.Range("H102").Value = -1200
$ws.Range("I102").Value = -2200
